$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell AD1: extend date series, preserving header style (s=1) ---
$ws.Range("AD1").Value = "'2020-04-03"
$ws.Range("AC1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Updated values in existing cells (recomputed cumulative counts) ---
$ws.Cells.Item(3, 21).Value = 7
$ws.Cells.Item(3, 22).Value = 8
$ws.Cells.Item(3, 23).Value = 8
$ws.Cells.Item(3, 24).Value = 8
$ws.Cells.Item(3, 25).Value = 8
$ws.Cells.Item(3, 26).Value = 9
$ws.Cells.Item(3, 27).Value = 9
$ws.Cells.Item(3, 28).Value = 10
$ws.Cells.Item(3, 29).Value = 11
$ws.Cells.Item(4, 20).Value = 3
$ws.Cells.Item(4, 21).Value = 3
$ws.Cells.Item(4, 22).Value = 4
$ws.Cells.Item(4, 23).Value = 4
$ws.Cells.Item(4, 24).Value = 4
$ws.Cells.Item(4, 25).Value = 6
$ws.Cells.Item(4, 26).Value = 8
$ws.Cells.Item(4, 27).Value = 8
$ws.Cells.Item(4, 28).Value = 11
$ws.Cells.Item(4, 29).Value = 12
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 3
$ws.Cells.Item(10, 19).Value = 3
$ws.Cells.Item(10, 20).Value = 5
$ws.Cells.Item(10, 21).Value = 5
$ws.Cells.Item(10, 22).Value = 7
$ws.Cells.Item(10, 23).Value = 7
$ws.Cells.Item(10, 24).Value = 7
$ws.Cells.Item(10, 25).Value = 8
$ws.Cells.Item(10, 26).Value = 8
$ws.Cells.Item(10, 27).Value = 9
$ws.Cells.Item(10, 28).Value = 9
$ws.Cells.Item(10, 29).Value = 10
$ws.Cells.Item(11, 16).Value = 8
$ws.Cells.Item(11, 17).Value = 8
$ws.Cells.Item(11, 18).Value = 8
$ws.Cells.Item(11, 19).Value = 9
$ws.Cells.Item(11, 20).Value = 9
$ws.Cells.Item(11, 21).Value = 9
$ws.Cells.Item(11, 22).Value = 10
$ws.Cells.Item(11, 23).Value = 13
$ws.Cells.Item(11, 24).Value = 14
$ws.Cells.Item(11, 25).Value = 14
$ws.Cells.Item(11, 26).Value = 14
$ws.Cells.Item(11, 27).Value = 16
$ws.Cells.Item(11, 28).Value = 18
$ws.Cells.Item(11, 29).Value = 19
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 3
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 3
$ws.Cells.Item(12, 7).Value = 3
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 4
$ws.Cells.Item(12, 10).Value = 6
$ws.Cells.Item(12, 11).Value = 11
$ws.Cells.Item(12, 12).Value = 13
$ws.Cells.Item(12, 13).Value = 15
$ws.Cells.Item(12, 14).Value = 23
$ws.Cells.Item(12, 15).Value = 26
$ws.Cells.Item(12, 16).Value = 36
$ws.Cells.Item(12, 17).Value = 38
$ws.Cells.Item(12, 18).Value = 45
$ws.Cells.Item(12, 19).Value = 65
$ws.Cells.Item(12, 20).Value = 78
$ws.Cells.Item(12, 21).Value = 91
$ws.Cells.Item(12, 22).Value = 107
$ws.Cells.Item(12, 23).Value = 119
$ws.Cells.Item(12, 24).Value = 133
$ws.Cells.Item(12, 25).Value = 147
$ws.Cells.Item(12, 26).Value = 163
$ws.Cells.Item(12, 27).Value = 176
$ws.Cells.Item(12, 28).Value = 191
$ws.Cells.Item(12, 29).Value = 211
$ws.Cells.Item(13, 16).Value = 4
$ws.Cells.Item(13, 17).Value = 7
$ws.Cells.Item(13, 18).Value = 10
$ws.Cells.Item(13, 19).Value = 12
$ws.Cells.Item(13, 20).Value = 18
$ws.Cells.Item(13, 21).Value = 20
$ws.Cells.Item(13, 22).Value = 22
$ws.Cells.Item(13, 23).Value = 24
$ws.Cells.Item(13, 24).Value = 27
$ws.Cells.Item(13, 25).Value = 30
$ws.Cells.Item(13, 26).Value = 33
$ws.Cells.Item(13, 27).Value = 37
$ws.Cells.Item(13, 28).Value = 40
$ws.Cells.Item(13, 29).Value = 42
$ws.Cells.Item(14, 21).Value = 8
$ws.Cells.Item(14, 22).Value = 9
$ws.Cells.Item(14, 23).Value = 9
$ws.Cells.Item(14, 24).Value = 9
$ws.Cells.Item(14, 25).Value = 11
$ws.Cells.Item(14, 26).Value = 14
$ws.Cells.Item(14, 27).Value = 15
$ws.Cells.Item(14, 28).Value = 16
$ws.Cells.Item(14, 29).Value = 20
$ws.Cells.Item(17, 20).Value = 2
$ws.Cells.Item(17, 21).Value = 2
$ws.Cells.Item(17, 22).Value = 2
$ws.Cells.Item(17, 23).Value = 2
$ws.Cells.Item(17, 24).Value = 2
$ws.Cells.Item(17, 25).Value = 3
$ws.Cells.Item(17, 26).Value = 3
$ws.Cells.Item(17, 27).Value = 3
$ws.Cells.Item(17, 28).Value = 4
$ws.Cells.Item(17, 29).Value = 5
$ws.Cells.Item(18, 25).Value = 10
$ws.Cells.Item(18, 26).Value = 10
$ws.Cells.Item(18, 27).Value = 10
$ws.Cells.Item(18, 28).Value = 12
$ws.Cells.Item(18, 29).Value = 12
$ws.Cells.Item(19, 21).Value = 6
$ws.Cells.Item(19, 22).Value = 7
$ws.Cells.Item(19, 23).Value = 7
$ws.Cells.Item(19, 24).Value = 7
$ws.Cells.Item(19, 25).Value = 9
$ws.Cells.Item(19, 26).Value = 9
$ws.Cells.Item(19, 27).Value = 9
$ws.Cells.Item(19, 28).Value = 10
$ws.Cells.Item(19, 29).Value = 12
$ws.Cells.Item(20, 25).Value = 27
$ws.Cells.Item(20, 26).Value = 27
$ws.Cells.Item(20, 27).Value = 28
$ws.Cells.Item(20, 28).Value = 28
$ws.Cells.Item(20, 29).Value = 34
$ws.Cells.Item(21, 15).Value = 8
$ws.Cells.Item(21, 16).Value = 9
$ws.Cells.Item(21, 17).Value = 12
$ws.Cells.Item(21, 18).Value = 14
$ws.Cells.Item(21, 19).Value = 18
$ws.Cells.Item(21, 20).Value = 19
$ws.Cells.Item(21, 21).Value = 20
$ws.Cells.Item(21, 22).Value = 25
$ws.Cells.Item(21, 23).Value = 27
$ws.Cells.Item(21, 24).Value = 31
$ws.Cells.Item(21, 25).Value = 37
$ws.Cells.Item(21, 26).Value = 40
$ws.Cells.Item(21, 27).Value = 44
$ws.Cells.Item(21, 28).Value = 49
$ws.Cells.Item(21, 29).Value = 50
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = 4
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = 6
$ws.Cells.Item(22, 7).Value = 6
$ws.Cells.Item(22, 8).Value = 7
$ws.Cells.Item(22, 9).Value = 9
$ws.Cells.Item(22, 10).Value = 14
$ws.Cells.Item(22, 11).Value = 20
$ws.Cells.Item(22, 12).Value = 25
$ws.Cells.Item(22, 13).Value = 29
$ws.Cells.Item(22, 14).Value = 44
$ws.Cells.Item(22, 15).Value = 52
$ws.Cells.Item(22, 16).Value = 71
$ws.Cells.Item(22, 17).Value = 83
$ws.Cells.Item(22, 18).Value = 104
$ws.Cells.Item(22, 19).Value = 142
$ws.Cells.Item(22, 20).Value = 176
$ws.Cells.Item(22, 21).Value = 206
$ws.Cells.Item(22, 22).Value = 245
$ws.Cells.Item(22, 23).Value = 272
$ws.Cells.Item(22, 24).Value = 300
$ws.Cells.Item(22, 25).Value = 336
$ws.Cells.Item(22, 26).Value = 364
$ws.Cells.Item(22, 27).Value = 391
$ws.Cells.Item(22, 28).Value = 427
$ws.Cells.Item(22, 29).Value = 477

# --- New column AD values (rows 2-23), new date column added ---
$ws.Cells.Item(2, 30).Value = 2
$ws.Cells.Item(3, 30).Value = 11
$ws.Cells.Item(4, 30).Value = 14
$ws.Cells.Item(5, 30).Value = 10
$ws.Cells.Item(6, 30).Value = 2
$ws.Cells.Item(7, 30).Value = 17
$ws.Cells.Item(8, 30).Value = 4
$ws.Cells.Item(9, 30).Value = 0
$ws.Cells.Item(10, 30).Value = 10
$ws.Cells.Item(11, 30).Value = 19
$ws.Cells.Item(12, 30).Value = 229
$ws.Cells.Item(13, 30).Value = 45
$ws.Cells.Item(14, 30).Value = 23
$ws.Cells.Item(15, 30).Value = 5
$ws.Cells.Item(16, 30).Value = 9
$ws.Cells.Item(17, 30).Value = 5
$ws.Cells.Item(18, 30).Value = 12
$ws.Cells.Item(19, 30).Value = 12
$ws.Cells.Item(20, 30).Value = 36
$ws.Cells.Item(21, 30).Value = 54
$ws.Cells.Item(22, 30).Value = 514
$ws.Cells.Item(23, 30).Value = 1
